$wb = $excel.ActiveWorkbook
$main = $wb.Worksheets.Item("Main")

# --- Update shared string "Cause" -> "Index Event" on N2 ---
$main.Range("N2").Value = "Index Event"

# --- Set N5 to 2 ---
$main.Range("N5").Value = 2

# --- Add rows 20-21 (Enable/Led on, Disable/Led off) ---
$main.Range("A20").Value = "Enable"
$main.Range("B20").Value = "Led on "
$main.Range("A21").Value = "Disable"
$main.Range("B21").Value = "Led off"

# --- Add the new EEPROM worksheet after Main ---
$eeprom = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $main)
$eeprom.Name = "EEPROM"

# Center-align the cells that will hold data (matches the rest of the workbook's style)
$eeprom.Range("A1:A10").HorizontalAlignment = -4108
$eeprom.Range("A1:A10").VerticalAlignment = -4108
$eeprom.Range("B1:B2").HorizontalAlignment = -4108
$eeprom.Range("B1:B2").VerticalAlignment = -4108

$eeprom.Range("A1").Value = "Address"
$eeprom.Range("B1").Value = "Variable"
$eeprom.Range("A2").Value = 0
$eeprom.Range("B2").Value = "state"
$eeprom.Range("A3").Value = 1
$eeprom.Range("A4").Value = 2
$eeprom.Range("A5").Value = 3
$eeprom.Range("A6").Value = 4
$eeprom.Range("A7").Value = 5
$eeprom.Range("A8").Value = 6
$eeprom.Range("A9").Value = 7
$eeprom.Range("A10").Value = 8

$eeprom.Range("B2").Select()

# --- Restore selection & active sheet on Main ---
$main.Range("B4:E18").Select()
$main.Activate()
